# edit.ps1 - apply LOB1223.xlsx content restructuring
# Rewrites the "Trabalho de Graduacao I" syllabus rows to insert the new
# Portuguese "Objetivos" paragraph, re-flow all downstream rows by one slot,
# add the new "Programa" / "Criterio" / "Bibliografia" content, and push the
# requirement list down by one row (adding the final LOQ4247 entry at row 51).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LOB1223")

# --- Update / insert the cell values that moved or changed text ---
$ws.Range("B10").Value = "O Trabalho de Graduação tem por objetivo a integração, o aprofundamento e a aplicação dos conhecimentos adquiridos ao longo do curso, preparando e desenvolvendo a capacidade do aluno para a realização de atividades, que fazem parte do perfil de atuação profissional do engenheiro ambiental. O Trabalho de Graduação I poderá ser desenvolvido seguindo dois modelos: Modelo Artigo ou Modelo Produto.No Modelo Artigo, o aluno deverá ser capaz de apresentar a ideia principal, especificamente delimitar a questão científica, e que deve ser suportada pelo estado da arte. O projeto precisa ainda indicar o método a ser empregado na investigação com etapas previstas e os resultados esperados.No Modelo Produto, o aluno deverá ser capaz de utilizar tecnologias digitais ou outras metodologias, para obtenção de produto aplicável ou tecnologia que utilize conteúdo da Engenharia Ambiental ou áreas correlatas."
$ws.Range("C10").Value = "O Trabalho de Graduação tem por objetivo a integração, o aprofundamento e a aplicação dos conhecimentos adquiridos ao longo do curso, preparando e desenvolvendo a capacidade do aluno para a realização de atividades, que fazem parte do perfil de atuação profissional do engenheiro ambiental. O Trabalho de Graduação I poderá ser desenvolvido seguindo dois modelos: Modelo Artigo ou Modelo Produto.No Modelo Artigo, o aluno deverá ser capaz de apresentar a ideia principal, especificamente delimitar a questão científica, e que deve ser suportada pelo estado da arte. O projeto precisa ainda indicar o método a ser empregado na investigação com etapas previstas e os resultados esperados.No Modelo Produto, o aluno deverá ser capaz de utilizar tecnologias digitais ou outras metodologias, para obtenção de produto aplicável ou tecnologia que utilize conteúdo da Engenharia Ambiental ou áreas correlatas."
$ws.Range("B13").Value = "4780627 - Ana Lucia Gabas Ferreira"
$ws.Range("C13").Value = "4780627 - Ana Lucia Gabas Ferreira"
$ws.Range("A14").Value = "Programa resumido:"
$ws.Range("B14").Value = "O aluno deverá escolher o modelo a ser seguido para o desenvolvimento do Trabalho de Graduação I, Artigo ou Produto. Diante da escolha do modelo, o aluno deverá desenvolver conteúdo científico ou produto aplicável, ambos relacionados à Engenharia Ambiental."
$ws.Range("C14").Value = "O aluno deverá escolher o modelo a ser seguido para o desenvolvimento do Trabalho de Graduação I, Artigo ou Produto. Diante da escolha do modelo, o aluno deverá desenvolver conteúdo científico ou produto aplicável, ambos relacionados à Engenharia Ambiental."
$ws.Range("A15").Value = "Short syllabus:"
$ws.Range("B15").Value = "The student must choose the model to be followed for the development of the Graduation Work I, Article or Product. When choosing the model, the student must develop scientific content or applicable product, both related to Environmental Engineering."
$ws.Range("C15").Value = "The student must choose the model to be followed for the development of the Graduation Work I, Article or Product. When choosing the model, the student must develop scientific content or applicable product, both related to Environmental Engineering."
$ws.Range("A16").Value = "Programa:"
$ws.Range("B16").Value = "Para ambos os modelos, o aluno deverá elaborar projeto cujo tema seja relacionado ao conteúdo programático do curso de Engenharia ambiental, podendo ser um tópico de interesse técnico ou científico ou ainda um estudo de caso. O programa da disciplina é constituído pelas seguintes etapas: 1. Definição e registro do tema 2. Aprovação do tema e do orientador 3. Acompanhamento mensal da evolução das atividades do aluno pelo orientador. 4. Entrega do projeto, 5. Avaliação e atribuição de notas por comissão de avaliadores."
$ws.Range("C16").Value = "Para ambos os modelos, o aluno deverá elaborar projeto cujo tema seja relacionado ao conteúdo programático do curso de Engenharia ambiental, podendo ser um tópico de interesse técnico ou científico ou ainda um estudo de caso. O programa da disciplina é constituído pelas seguintes etapas: 1. Definição e registro do tema 2. Aprovação do tema e do orientador 3. Acompanhamento mensal da evolução das atividades do aluno pelo orientador. 4. Entrega do projeto, 5. Avaliação e atribuição de notas por comissão de avaliadores."
$ws.Range("A17").Value = "Syllabus:"
$ws.Range("B17").Value = "For both models, the student must prepare a project whose theme is related to the programmatic content of the Environmental Engineering course, which may be a topic of technical or scientific interest or even a case study. The course program consists of the following stages: 1. Definition and registration of the theme 2. Approval of the theme and the advisor 3. Monthly monitoring of the evolution of the student's activities by the advisor. 4. Project delivery, 5. Evaluation and grading by an evaluation committee."
$ws.Range("C17").Value = "For both models, the student must prepare a project whose theme is related to the programmatic content of the Environmental Engineering course, which may be a topic of technical or scientific interest or even a case study. The course program consists of the following stages: 1. Definition and registration of the theme 2. Approval of the theme and the advisor 3. Monthly monitoring of the evolution of the student's activities by the advisor. 4. Project delivery, 5. Evaluation and grading by an evaluation committee."
$ws.Range("A18").Value = "Avaliação:"
$ws.Range("A19").Value = "Método:"
$ws.Range("A20").Value = "Critério:"
$ws.Range("A21").Value = "Norma de recuperação:"
$ws.Range("A22").Value = "Bibliografia:"
$ws.Range("B22").Value = "A ser definida no plano de trabalho."
$ws.Range("C22").Value = "A ser definida no plano de trabalho."
$ws.Range("A23").Value = "Requisitos:"
$ws.Range("B24").Value = "LOB1003 -  Cálculo I  (Requisito fraco)`n"
$ws.Range("C24").Value = "LOB1003 -  Cálculo I  (Requisito fraco)`n"
$ws.Range("B25").Value = "LOB1004 -  Cálculo II  (Requisito fraco)`n"
$ws.Range("C25").Value = "LOB1004 -  Cálculo II  (Requisito fraco)`n"
$ws.Range("B26").Value = "LOB1006 -  Cálculo IV  (Requisito fraco)`n"
$ws.Range("C26").Value = "LOB1006 -  Cálculo IV  (Requisito fraco)`n"
$ws.Range("B27").Value = "LOB1011 -  Eletricidade Aplicada  (Requisito fraco)`n"
$ws.Range("C27").Value = "LOB1011 -  Eletricidade Aplicada  (Requisito fraco)`n"
$ws.Range("B28").Value = "LOB1012 -  Estatística  (Requisito fraco)`n"
$ws.Range("C28").Value = "LOB1012 -  Estatística  (Requisito fraco)`n"
$ws.Range("B29").Value = "LOB1018 -  Física I  (Requisito fraco)`n"
$ws.Range("C29").Value = "LOB1018 -  Física I  (Requisito fraco)`n"
$ws.Range("B30").Value = "LOB1019 -  Física II  (Requisito fraco)`n"
$ws.Range("C30").Value = "LOB1019 -  Física II  (Requisito fraco)`n"
$ws.Range("B31").Value = "LOB1021 -  Física IV  (Requisito fraco)`n"
$ws.Range("C31").Value = "LOB1021 -  Física IV  (Requisito fraco)`n"
$ws.Range("B32").Value = "LOB1024 -  Mecânica  (Requisito fraco)`n"
$ws.Range("C32").Value = "LOB1024 -  Mecânica  (Requisito fraco)`n"
$ws.Range("B33").Value = "LOB1036 -  Geometria Analítica  (Requisito fraco)`n"
$ws.Range("C33").Value = "LOB1036 -  Geometria Analítica  (Requisito fraco)`n"
$ws.Range("B34").Value = "LOB1037 -  Àlgebra Linear  (Requisito fraco)`n"
$ws.Range("C34").Value = "LOB1037 -  Àlgebra Linear  (Requisito fraco)`n"
$ws.Range("B35").Value = "LOB1038 -  Física Experimental I  (Requisito fraco)`n"
$ws.Range("C35").Value = "LOB1038 -  Física Experimental I  (Requisito fraco)`n"
$ws.Range("B36").Value = "LOB1039 -  Física Experimental III  (Requisito fraco)`n"
$ws.Range("C36").Value = "LOB1039 -  Física Experimental III  (Requisito fraco)`n"
$ws.Range("B37").Value = "LOB1041 -  Física Experimental II  (Requisito fraco)`n"
$ws.Range("C37").Value = "LOB1041 -  Física Experimental II  (Requisito fraco)`n"
$ws.Range("B38").Value = "LOB1042 -  Física Experimental IV  (Requisito fraco)`n"
$ws.Range("C38").Value = "LOB1042 -  Física Experimental IV  (Requisito fraco)`n"
$ws.Range("B39").Value = "LOB1045 -  Leitura e Produção de Textos Acadêmicos  (Requisito fraco)`n"
$ws.Range("C39").Value = "LOB1045 -  Leitura e Produção de Textos Acadêmicos  (Requisito fraco)`n"
$ws.Range("B40").Value = "LOB1052 -  Cálculo III  (Requisito fraco)`n"
$ws.Range("C40").Value = "LOB1052 -  Cálculo III  (Requisito fraco)`n"
$ws.Range("B41").Value = "LOB1053 -  Física III  (Requisito fraco)`n"
$ws.Range("C41").Value = "LOB1053 -  Física III  (Requisito fraco)`n"
$ws.Range("B42").Value = "LOB1056 -  Introdução aos Métodos Numéricos e Computacionais  (Requisito fraco)`n"
$ws.Range("C42").Value = "LOB1056 -  Introdução aos Métodos Numéricos e Computacionais  (Requisito fraco)`n"
$ws.Range("B43").Value = "LOB1232 -  Licenciamento Ambiental  (Requisito fraco)`n"
$ws.Range("C43").Value = "LOB1232 -  Licenciamento Ambiental  (Requisito fraco)`n"
$ws.Range("B44").Value = "LOB1257 -  Sistema de Abastecimento e Tratamento de Água  (Requisito fraco)`n"
$ws.Range("C44").Value = "LOB1257 -  Sistema de Abastecimento e Tratamento de Água  (Requisito fraco)`n"
$ws.Range("B45").Value = "LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito fraco)`n"
$ws.Range("C45").Value = "LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito fraco)`n"
$ws.Range("B46").Value = "LOM3081 -  Introdução à Mecânica dos Sólidos  (Requisito fraco)`n"
$ws.Range("C46").Value = "LOM3081 -  Introdução à Mecânica dos Sólidos  (Requisito fraco)`n"
$ws.Range("B47").Value = "LOQ4095 -  Química Geral Experimental  (Requisito fraco)`n"
$ws.Range("C47").Value = "LOQ4095 -  Química Geral Experimental  (Requisito fraco)`n"
$ws.Range("B48").Value = "LOQ4097 -  Fundamentos de Química para Engenharia I (Requisito fraco)`n"
$ws.Range("C48").Value = "LOQ4097 -  Fundamentos de Química para Engenharia I (Requisito fraco)`n"
$ws.Range("B49").Value = "LOQ4098 -  Fundamentos de Química para Engenharia II (Requisito fraco)`n"
$ws.Range("C49").Value = "LOQ4098 -  Fundamentos de Química para Engenharia II (Requisito fraco)`n"
$ws.Range("B50").Value = "LOQ4233 -  Gestão de Negócios  (Requisito fraco)`n"
$ws.Range("C50").Value = "LOQ4233 -  Gestão de Negócios  (Requisito fraco)`n"
$ws.Range("B51").Value = "LOQ4247 -  Desenho Assistido por Computador  (Requisito fraco)`n"
$ws.Range("C51").Value = "LOQ4247 -  Desenho Assistido por Computador  (Requisito fraco)`n"

# --- Clear cells that become empty after the reflow ---
$ws.Range("A13").Clear()
$ws.Range("B18").Clear()
$ws.Range("C18").Clear()
$ws.Range("B23").Clear()
$ws.Range("C23").Clear()

# --- Row height adjustments to match the new layout ---
$ws.Rows.Item(13).AutoFit()
$ws.Rows.Item(15).RowHeight = 60
$ws.Rows.Item(17).RowHeight = 120
$ws.Rows.Item(18).AutoFit()
$ws.Rows.Item(21).RowHeight = 60
$ws.Rows.Item(22).RowHeight = 120
$ws.Rows.Item(23).AutoFit()
$ws.Rows.Item(51).RowHeight = 30
